$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (monto_invertido), shifting it (and
# monto_vencimiento) one column to the right so the new layout becomes:
# A: estado, B: moneda, C: orden, D: monto_invertido, E: monto_vencimiento
$ws.Columns("C:C").Insert()

# Match the width Excel would carry over to the newly-inserted column from
# its left-hand neighbour (same default width as columns A/B).
$ws.Columns("C:C").ColumnWidth = $ws.Columns("A:A").ColumnWidth

# Header for the new column - same style as the other plain headers (A1/B1).
$ws.Range("C1").Value = "orden"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").WrapText = $true

# Fill in the order numbers 1..8 for the data rows, using the same (non-bold)
# style as the other data columns (A2:B9).
for ($i = 2; $i -le 9; $i++) {
    $ws.Cells.Item($i, 3).Value = $i - 1
}
$ws.Range("C2:C9").Font.Bold = $false
$ws.Range("C2:C9").WrapText = $true

# Update selection to mirror the saved workbook state.
$ws.Range("E11").Select()
